# "Colocando header nos gráficos"
# Add a header label in column A (row 1) for each data table, fix a few
# accented Portuguese words, remove the now-redundant bold/border style
# from the row-label cells (A2:A12, etc. - only the real header row keeps
# the style), and apply the small data corrections described in the diff.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheets 1-4 share an identical layout:
#   A1 (new header) | B1..E1 = years
#   A2:A12 = source/technology names
# ---------------------------------------------------------------------
$sheetIndexes = @(1, 2, 3, 4)

foreach ($idx in $sheetIndexes) {
    $ws = $wb.Worksheets.Item($idx)

    # Add the new column header in A1, copying the existing header format
    # from B1 so it matches the other header cells exactly.
    $ws.Range("B1").Copy()
    $ws.Range("A1").PasteSpecial(-4122)
    $ws.Range("A1").Value = "Fonte/Tecnologia"

    # Fix accentuation on a few labels.
    $ws.Range("A3").Value = "Gás Natural"
    $ws.Range("A4").Value = "Carvão"
    $ws.Range("A6").Value = "Óleos Comb"
    $ws.Range("A8").Value = "Eólica"
    $ws.Range("A11").Value = "Pot. Compl."

    # The row-label cells are no longer the "header" row, so drop their
    # bold/border formatting (A1 is now the only header in column A).
    $ws.Range("A2:A12").ClearFormats()
}

# ---------------------------------------------------------------------
# Sheet 5: "Emissoes Totais (MtCO2eq)"
# ---------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item(5)

$ws5.Range("B1").Copy()
$ws5.Range("A1").PasteSpecial(-4122)
$ws5.Range("A1").Value = "Período"

$ws5.Range("A2").Value = "P.Médio"
$ws5.Range("A3").Value = "P.Crítico"
$ws5.Range("A2:A3").ClearFormats()

# Remove the obsolete "Teto" row entirely.
$ws5.Rows.Item(4).Delete()

# ---------------------------------------------------------------------
# Sheet 6: "Custo Total (bilhões de R$)"
# ---------------------------------------------------------------------
$ws6 = $wb.Worksheets.Item(6)

$ws6.Range("B1").Copy()
$ws6.Range("A1").PasteSpecial(-4122)
$ws6.Range("A1").Value = "Tipo Expansão"

# B1's label changes from "Custo" to "2015" - keep it a text value (like
# the "2015" header used on the other sheets) rather than a number, and
# keep the existing header style untouched.
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B1").Copy()
$ws6.Range("B1").PasteSpecial(-4163)

$ws6.Range("A2").Value = "Expansão Centralizada"
$ws6.Range("B2").Value = 743
$ws6.Range("A3").Value = "Expansão por GD"
$ws6.Range("B3").Value = 99
$ws6.Range("A2:A3").ClearFormats()
